# Moved statistical datasets and results
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A (shifts existing data A:D -> B:E)
$ws.Range("A1").EntireColumn.Insert()

# Insert a new row before row 1 (shifts existing data down by one row)
$ws.Range("A1").EntireRow.Insert()

# New header row
$ws.Range("B1").Value = "Valid"
$ws.Range("C1").Value = "T"
$ws.Range("D1").Value = "Z"
$ws.Range("E1").Value = "p-value"

# New label column (column A, rows 2-16)
$labels = @(
    "CyclomaticComplexity(CC) & NbOperators",
    "CyclomaticComplexity(CC) & EffortToImplement",
    "MaintainabilityIndex & MaintainabilityIndex",
    "NbOperands & NbOperands",
    "NbOperands & EffortToImplement",
    "NbOperators & CyclomaticComplexity(CC)",
    "NbOperators & EffortToImplement",
    "DifficultyLevel & DifficultyLevel",
    "ProgramLevel & ProgramLevel",
    "EffortToImplement & CyclomaticComplexity(CC)",
    "EffortToImplement & NbOperands",
    "EffortToImplement & NbOperators",
    "EffortToImplement & ProgramLength",
    "EffortToImplement & EffortToImplement",
    "TimeToImplement & TimeToImplement"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $labels[$i]
}

# Set the width of the newly inserted label column (columns B:E already retain
# their original widths because they were shifted, not re-created, by the insert)
$ws.Columns.Item(1).ColumnWidth = 53.7
